$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11, columns A-G
$data = @(
    @(97105, "Dra. Olivia Nogueira", "Engenharia",        "Viagem de negócios", 2, 45086, 2815.52),
    @(92017, "Igor Oliveira",        "Operações",         "Doença",             3, 45105, 4965.95),
    @(72675, "Maria Julia Cardoso",  "Recursos Humanos",  "Problemas pessoais", 7, 45100, 6777.92),
    @(18336, "Dr. Davi Silveira",    "P&D",                "Consulta médica",    8, 45085, 3313.77),
    @(63121, "Kevin Pinto",          "Engenharia",        "Viagem de negócios", 1, 45079, 5172.9),
    @(24709, "Davi Luiz Campos",     "Financeiro",        "Doença",             3, 45100, 10941.97),
    @(74219, "Maria Vitória Pires",  "Operações",         "Viagem de negócios", 2, 45085, 7637.73),
    @(87149, "Ana Clara Pereira",    "Engenharia",        "Problemas pessoais", 4, 45103, 4870.72),
    @(82721, "Sr. João Melo",        "TI",                "Problemas pessoais", 4, 45081, 5259.45),
    @(42236, "Ana Júlia Barros",     "Vendas",             "Consulta médica",    2, 45098, 5981.86)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
